$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 276
$ws.Range("C3").Value = 173083
$ws.Range("C4").Value = 163790
$ws.Range("C5").Value = 9293
$ws.Range("C6").Value = 510
$ws.Range("C7").Value = 5.37
$ws.Range("C8").Value = 65.94
